$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103 (pushes existing row 103 and below down by one)
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 - it's a duplicate of the (old) row 103 data
# except for the Fecha (date) column D, which gets a new date serial value.
$ws.Cells.Item(103, 1).Value = 7
$ws.Cells.Item(103, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(103, 3).Value = "Ñuble"
$ws.Cells.Item(103, 4).Value = 44452
$ws.Cells.Item(103, 5).Value = 16
$ws.Cells.Item(103, 6).Value = 100112023
$ws.Cells.Item(103, 7).Value = "Brócoli"
$ws.Cells.Item(103, 8).Value = "Sin especificar"
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value = 300
$ws.Cells.Item(103, 11).Value = 600
$ws.Cells.Item(103, 12).Value = 650
$ws.Cells.Item(103, 13).Value = 625
$ws.Cells.Item(103, 14).Value = "`$/unidad"
$ws.Cells.Item(103, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(103, 16).Value = 625
$ws.Cells.Item(103, 17).Value = 1
$ws.Cells.Item(103, 18).Value = "Hortaliza"
